# Natmi following Dr Hou advice
# Updates recalculated NATMI LR-pair statistics for Hspg2-Ptprs (OldD7) after
# re-running the pipeline with the advised parameters (ligand/receptor-expressing
# cell counts change from 1 to 3 per cluster, with all derived columns updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "E"=3; "G"=155.2138263333333; "H"=465.641479; "I"=0.3492508712612995; "J"=0.3492508712612995; "K"=3; "M"=4.127188333333333; "N"=12.381565; "O"=0.0561359176022362; "P"=0.05613591760223619; "Q"=640.5966932149594; "R"=5765.370238934634; "S"=0.01960551813163351; "T"=0.01960551813163351 }
    3 = @{ "E"=3; "G"=155.2138263333333; "H"=465.641479; "I"=0.3492508712612995; "J"=0.3492508712612995; "K"=3; "M"=47.24901333333333; "N"=141.74704; "O"=0.6426570597336346; "P"=0.6426570597336345; "Q"=7333.700149941351; "R"=66003.30134947215; "S"=0.2244485380341969; "T"=0.2244485380341968 }
    4 = @{ "E"=3; "G"=155.2138263333333; "H"=465.641479; "I"=0.3492508712612995; "J"=0.3492508712612995; "K"=3; "M"=8.218847999999999; "N"=24.656544; "O"=0.1117885923419141; "P"=0.1117885923419141; "Q"=1275.678846132064; "R"=11481.10961518858; "S"=0.03904226327248773; "T"=0.03904226327248773 }
    5 = @{ "E"=3; "G"=155.2138263333333; "H"=465.641479; "I"=0.3492508712612995; "J"=0.3492508712612995; "K"=3; "M"=13.92629833333334; "N"=41.77889500000001; "O"=0.1894184303222152; "P"=0.1894184303222152; "Q"=2161.55405097619; "R"=19453.98645878571; "S"=0.06615455182298143; "T"=0.06615455182298141 }
    6 = @{ "E"=3; "G"=246.1811623333333; "H"=738.543487; "I"=0.5539389593320749; "J"=0.5539389593320749; "K"=3; "M"=4.127188333333333; "N"=12.381565; "O"=0.0561359176022362; "P"=0.05613591760223619; "Q"=1016.036021068573; "R"=9144.324189617155; "S"=0.03109587177773383; "T"=0.03109587177773382 }
    7 = @{ "E"=3; "G"=246.1811623333333; "H"=738.543487; "I"=0.5539389593320749; "J"=0.5539389593320749; "K"=3; "M"=47.24901333333333; "N"=141.74704; "O"=0.6426570597336346; "P"=0.6426570597336345; "Q"=11631.81702150316; "R"=104686.3531935285; "S"=0.3559927828762607; "T"=0.3559927828762606 }
    8 = @{ "E"=3; "G"=246.1811623333333; "H"=738.543487; "I"=0.5539389593320749; "J"=0.5539389593320749; "K"=3; "M"=8.218847999999999; "N"=24.656544; "O"=0.1117885923419141; "P"=0.1117885923419141; "Q"=2023.325553680992; "R"=18209.92998312893; "S"=0.06192405650707745; "T"=0.06192405650707745 }
    9 = @{ "E"=3; "G"=246.1811623333333; "H"=738.543487; "I"=0.5539389593320749; "J"=0.5539389593320749; "K"=3; "M"=13.92629833333334; "N"=41.77889500000001; "O"=0.1894184303222152; "P"=0.1894184303222152; "Q"=3428.392310700763; "R"=30855.53079630687; "S"=0.104926248171003; "T"=0.104926248171003 }
    10 = @{ "E"=3; "G"=0.2401933333333333; "H"=0.72058; "I"=0.0005404655817044752; "J"=0.0005404655817044752; "K"=3; "M"=4.127188333333333; "N"=12.381565; "O"=0.0561359176022362; "P"=0.05613591760223619; "Q"=0.9913231230777777; "R"=8.921908107699998; "S"=0.00003033953136140708; "T"=0.00003033953136140707 }
    11 = @{ "E"=3; "G"=0.2401933333333333; "H"=0.72058; "I"=0.0005404655817044752; "J"=0.0005404655817044752; "K"=3; "M"=47.24901333333333; "N"=141.74704; "O"=0.6426570597336346; "P"=0.6426570597336345; "Q"=11.34889800924444; "R"=102.1400820832; "S"=0.0003473340216254265; "T"=0.0003473340216254265 }
    12 = @{ "E"=3; "G"=0.2401933333333333; "H"=0.72058; "I"=0.0005404655817044752; "J"=0.0005404655817044752; "K"=3; "M"=8.218847999999999; "N"=24.656544; "O"=0.1117885923419141; "P"=0.1117885923419141; "Q"=1.97411249728; "R"=17.76701247552; "S"=0.00006041788658799704; "T"=0.00006041788658799704 }
    13 = @{ "E"=3; "G"=0.2401933333333333; "H"=0.72058; "I"=0.0005404655817044752; "J"=0.0005404655817044752; "K"=3; "M"=13.92629833333334; "N"=41.77889500000001; "O"=0.1894184303222152; "P"=0.1894184303222152; "Q"=3.345004017677778; "R"=30.1050361591; "S"=0.0001023741421296447; "T"=0.0001023741421296446 }
    14 = @{ "E"=3; "G"=42.784114; "H"=128.352342; "I"=0.09626970382492123; "J"=0.09626970382492124; "K"=3; "M"=4.127188333333333; "N"=12.381565; "O"=0.0561359176022362; "P"=0.05613591760223619; "Q"=176.5780961528033; "R"=1589.20286537523; "S"=0.005404188161507461; "T"=0.005404188161507461 }
    15 = @{ "E"=3; "G"=42.784114; "H"=128.352342; "I"=0.09626970382492123; "J"=0.09626970382492124; "K"=3; "M"=47.24901333333333; "N"=141.74704; "O"=0.6426570597336346; "P"=0.6426570597336345; "Q"=2021.507172840853; "R"=18193.56455556768; "S"=0.06186840480155172; "T"=0.06186840480155172 }
    16 = @{ "E"=3; "G"=42.784114; "H"=128.352342; "I"=0.09626970382492123; "J"=0.09626970382492124; "K"=3; "M"=8.218847999999999; "N"=24.656544; "O"=0.1117885923419141; "P"=0.1117885923419141; "Q"=351.6361297806719; "R"=3164.725168026048; "S"=0.01076185467576093; "T"=0.01076185467576093 }
    17 = @{ "E"=3; "G"=42.784114; "H"=128.352342; "I"=0.09626970382492123; "J"=0.09626970382492124; "K"=3; "M"=13.92629833333334; "N"=41.77889500000001; "O"=0.1894184303222152; "P"=0.1894184303222152; "Q"=595.8243354913434; "R"=5362.41901942209; "S"=0.01823525618610114; "T"=0.01823525618610114 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
